$d = $word.ActiveDocument

# Near the end of the document there is a run of identical empty
# paragraphs (rFonts majorHAnsi/Times New Roman, szCs=24) immediately
# followed by one final paragraph that contains the "contact us"
# sentence with wngkai.91@gmail.com. The edit:
#   1) removes 11 of those empty paragraphs, and
#   2) clears the e-mail sentence out of the final paragraph, leaving
#      the (now empty) paragraph and its formatting in place.

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# --- Step 1: delete the 11 empty paragraphs immediately preceding the last one.
$endPos = $lastPara.Range.Start
$startPara = $d.Paragraphs.Item($count - 11)
$startPos = $startPara.Range.Start
$d.Range($startPos, $endPos).Delete()

# --- Step 2: remove the e-mail sentence text/run from the final paragraph.
$lastParaNow = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastParaNow.Range.Find.Execute(
    "If you are having difficulties viewing our library database system, please contact wngkai.91@gmail.com",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
